$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a new value that looks like a plain number must be forced to
# Text format first, otherwise Excel will silently convert the literal into a
# numeric cell instead of keeping it as the original text value.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.548.39"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.881.19"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "246.41"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4725"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.2887"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "0.06528"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "22.11"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.7725"
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "100.74"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("D13").Value = "0.07830"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "1.878.42"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "5.241"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "285.08"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "30.545.82"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "13.20"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "0.000007526"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "2.122.62"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "5.372"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "6.400"
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").Value = "9.140"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").Value = "162.74"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "1.917"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "0.09705"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "1.504"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "4.267"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "4.198"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "0.04846"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "0.6961"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "2.754"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "2.874"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").Value = "76.36"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "6.278"
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "1.980"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "0.4263"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "0.8313"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "101.55"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "9.786"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "7.044"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "35.11"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "891.27"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").Value = "0.05764"
$ws.Range("E51").Value = "  +0.17%  "

# Restore the default (unstyled) look for the cells we temporarily switched to
# Text format, now that the literal text value has been safely stored.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
